$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings (single decimal point) are kept as
# literal text (matching the source data feed formatting) instead of being
# auto-converted to floating point numbers by Excel.

$ws.Range('D2').Value = '30.417.71'
$ws.Range('E2').Value = '  -1.11%  '

$ws.Range('D3').Value = '1.916.11'
$ws.Range('E3').Value = '  +1.65%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9996'
$ws.Range('E4').Value = '  +0.05%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '242.34'
$ws.Range('E5').Value = '  +1.46%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9998'

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4697'
$ws.Range('E7').Value = '  -1.32%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2855'
$ws.Range('E8').Value = '  -0.41%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06811'
$ws.Range('E9').Value = '  +3.83%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '107.77'
$ws.Range('E10').Value = '  +12.25%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '18.35'
$ws.Range('E11').Value = '  -2.77%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07717'
$ws.Range('E12').Value = '  +1.97%  '

$ws.Range('D13').Value = '1.888.24'
$ws.Range('E13').Value = '  +0.23%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.244'
$ws.Range('E14').Value = '  +2.61%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6582'
$ws.Range('E15').Value = '  +0.57%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '295.28'
$ws.Range('E16').Value = '  -3.64%  '

$ws.Range('D17').Value = '30.440.07'
$ws.Range('E17').Value = '  -1.03%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000007616'
$ws.Range('E18').Value = '  +0.51%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.9994'
$ws.Range('E19').Value = '  +0.00%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.92'
$ws.Range('E20').Value = '  -1.92%  '

$ws.Range('D21').Value = '2.139.62'
$ws.Range('E21').Value = '  +0.52%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9999'
$ws.Range('E22').Value = '  +0.14%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.236'
$ws.Range('E23').Value = '  +2.26%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.241'
$ws.Range('E24').Value = '  +1.19%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.297'
$ws.Range('E25').Value = '  +0.22%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '167.96'
$ws.Range('E26').Value = '  +0.78%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '21.52'
$ws.Range('E27').Value = '  +6.70%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.084'
$ws.Range('E28').Value = '  +6.97%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.1070'
$ws.Range('E29').Value = '  -0.52%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.368'

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.158'
$ws.Range('E31').Value = '  -0.15%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.978'
$ws.Range('E32').Value = '  +0.39%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05067'

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7383'
$ws.Range('E34').Value = '  +1.05%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.152'
$ws.Range('E35').Value = '  -1.75%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02105'
$ws.Range('E36').Value = '  +8.73%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.747'
$ws.Range('E37').Value = '  +1.35%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.677'
$ws.Range('E38').Value = '  -0.87%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.058'
$ws.Range('E39').Value = '  -0.76%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '109.99'
$ws.Range('E40').Value = '  +2.17%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8722'
$ws.Range('E41').Value = '  -3.34%  '

$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.820'
$ws.Range('E42').Value = '  +3.34%  '

$ws.Range('B43').Value = 'TheSandbox'
$ws.Range('C43').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.4275'
$ws.Range('E43').Value = '  +1.41%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.9994'
$ws.Range('E44').Value = '  +0.01%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '67.44'
$ws.Range('E45').Value = '  +2.46%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '50.96'
$ws.Range('E46').Value = '  +18.66%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.180'
$ws.Range('E47').Value = '  -2.36%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.314'
$ws.Range('E48').Value = '  +3.26%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '34.88'
$ws.Range('E49').Value = '  +0.37%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.1209'
$ws.Range('E50').Value = '  -1.22%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.3943'
$ws.Range('E51').Value = '  +2.79%  '
